$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$__style = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.800.36'
$ws.Range('D2').Style = $__style
$ws.Range('E2').Value = '  -0.28%  '

$__style = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.046.45'
$ws.Range('D3').Style = $__style
$ws.Range('E3').Value = '  -1.26%  '

$ws.Range('E4').Value = '  +0.41%  '

$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '533.05'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  -1.70%  '

$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.84'
$ws.Range('D6').Style = $__style
$ws.Range('E6').Value = '  -0.19%  '

$ws.Range('E7').Value = '  +0.30%  '

$__style = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.042.74'
$ws.Range('D8').Style = $__style
$ws.Range('E8').Value = '  -1.27%  '

$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.489'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  -0.49%  '

$__style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('D10').Style = $__style
$ws.Range('E10').Value = '  -1.26%  '

$__style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.15'
$ws.Range('D11').Style = $__style
$ws.Range('E11').Value = '  -1.30%  '

$__style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('D12').Style = $__style
$ws.Range('E12').Value = '  -3.58%  '

$__style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000218'
$ws.Range('D13').Style = $__style
$ws.Range('E13').Value = '  -1.89%  '

$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.94'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  -3.08%  '

$__style = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.548.15'
$ws.Range('D15').Style = $__style
$ws.Range('E15').Value = '  +0.56%  '

$__style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.859.12'
$ws.Range('D16').Style = $__style
$ws.Range('E16').Value = '  +0.15%  '

$__style = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.112'
$ws.Range('D17').Style = $__style
$ws.Range('E17').Value = '  +1.14%  '

$__style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.060.19'
$ws.Range('D18').Style = $__style
$ws.Range('E18').Value = '  +0.00%  '

$__style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.55'
$ws.Range('D19').Style = $__style
$ws.Range('E19').Value = '  -1.93%  '

$__style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '465.20'
$ws.Range('D20').Style = $__style
$ws.Range('E20').Value = '  -4.05%  '

$__style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.24'
$ws.Range('D21').Style = $__style
$ws.Range('E21').Value = '  -1.69%  '

$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.685'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  -3.46%  '

$__style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.91'
$ws.Range('D23').Style = $__style
$ws.Range('E23').Value = '  -4.99%  '

$__style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.79'
$ws.Range('D24').Style = $__style
$ws.Range('E24').Value = '  -1.37%  '

$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.97'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  -1.59%  '

$ws.Range('E26').Value = '  -0.41%  '

$__style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.66'
$ws.Range('D27').Style = $__style
$ws.Range('E27').Value = '  -2.35%  '

$__style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.74'
$ws.Range('D28').Style = $__style
$ws.Range('E28').Value = '  -6.34%  '

$__style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = $__style
$ws.Range('E29').Value = '  +0.46%  '

$__style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '25.79'
$ws.Range('D30').Style = $__style
$ws.Range('E30').Value = '  -1.51%  '

$__style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.15'
$ws.Range('D31').Style = $__style
$ws.Range('E31').Value = '  +3.83%  '

$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.85'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  -4.51%  '

$__style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '58.20'
$ws.Range('D33').Style = $__style
$ws.Range('E33').Value = '  -1.37%  '

$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.27'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  -7.01%  '

$__style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.36'
$ws.Range('D35').Style = $__style
$ws.Range('E35').Value = '  +2.86%  '

$__style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.86'
$ws.Range('D36').Style = $__style
$ws.Range('E36').Value = '  -3.15%  '

$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '464.31'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  -3.37%  '

$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.214.20'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  +2.28%  '

$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0391'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  -0.48%  '

$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0782'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  -2.32%  '

$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.116'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  -0.15%  '

$__style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.02'
$ws.Range('D42').Style = $__style
$ws.Range('E42').Value = '  -1.11%  '

$__style = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.53'
$ws.Range('D43').Style = $__style
$ws.Range('E43').Value = '  -2.28%  '

$ws.Range('E44').Value = '  +0.13%  '

$__style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.247'
$ws.Range('D45').Style = $__style
$ws.Range('E45').Value = '  -3.04%  '

$__style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '122.69'
$ws.Range('D46').Style = $__style
$ws.Range('E46').Value = '  +3.30%  '

$__style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.75'
$ws.Range('D47').Style = $__style
$ws.Range('E47').Value = '  -0.64%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.108'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  -0.17%  '

$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$__style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.97'
$ws.Range('D49').Style = $__style
$ws.Range('E49').Value = '  -4.66%  '

$__style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0510'
$ws.Range('D50').Style = $__style
$ws.Range('E50').Value = '  -0.64%  '

$ws.Range('E51').Value = '  +4.66%  '
